# Update PvsI (Photosynthesis vs Irradiance) model-fitted respirometry rates.
# The "volume" (T), "rate.abs" (Z), "rate.a.spec" (AB) and "rate.output" (AD)
# columns are recalculated for rows 2-17 (Sheet1) to reflect the refitted model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ T = 0.1405756097560976;  Z = -0.1221848774731355;   AB = -492.9841458345362;  AD = -492.9841458345362 }
    3  = @{ T = 0.144009756097561;   Z = -0.1859906850082795;   AB = -1225.261743171314;  AD = -1225.261743171314 }
    4  = @{ T = 0.1462634146341464;  Z = -0.1176290256075774;   AB = -753.5407735840416;  AD = -753.5407735840416 }
    5  = @{ T = 0.1500390243902439;  Z = -0.1352576752579715;   AB = -844.7203620226234;  AD = -844.7203620226234 }
    6  = @{ T = 0.1463707317073171;  Z = -0.09541852482610373;  AB = -596.0777440014887;  AD = -596.0777440014887 }
    7  = @{ T = 0.1465658536585366;  Z = -0.1408857642960295;   AB = -703.0082285039052;  AD = -703.0082285039052 }
    8  = @{ T = 0.1449658536585366;  Z = -0.1744999617303195;   AB = -703.8747459606614;  AD = -703.8747459606614 }
    9  = @{ T = 0.1544;              Z = -0.0007693950938519304 }
    10 = @{ T = 0.1405756097560976;  Z = 0.1643696141026553;    AB = 663.1885670742363;   AD = 663.1885670742363 }
    11 = @{ T = 0.144009756097561;   Z = 0.121375857644185;     AB = 799.5948555671325;   AD = 799.5948555671325 }
    12 = @{ T = 0.1462634146341464;  Z = 0.07494526673124569;   AB = 480.1052629435983;   AD = 480.1052629435983 }
    13 = @{ T = 0.1500390243902439;  Z = 0.1314022337440853;    AB = 820.6420984774643;   AD = 820.6420984774643 }
    14 = @{ T = 0.1463707317073171;  Z = 0.08973606109305092;   AB = 560.5794991004969;   AD = 560.5794991004969 }
    15 = @{ T = 0.1465658536585366;  Z = 0.1444903154080865;    AB = 720.9946383054976;   AD = 720.9946383054976 }
    16 = @{ T = 0.1449658536585366;  Z = 0.2158516139991994;    AB = 870.6735432050632;   AD = 870.6735432050632 }
    17 = @{ T = 0.1544;              Z = -0.001158693818402261 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("T$row").Value = $vals.T
    $ws.Range("Z$row").Value = $vals.Z
    if ($vals.ContainsKey("AB")) {
        $ws.Range("AB$row").Value = $vals.AB
    }
    if ($vals.ContainsKey("AD")) {
        $ws.Range("AD$row").Value = $vals.AD
    }
}
